$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Geraldines GFC, Dundalk", "Louth"),
    @("Manguard Park Pitch 1", "Kildare"),
    @("Ballinascreen", "Derry"),
    @("Baltinglass", "Wicklow")
)

$row = 109
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
